$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2014-10"
$ws.Range("B2").Value = 100.037
$ws.Range("C2").Value = 99.3776
$ws.Range("D2").Value = 101.8211
$ws.Range("E2").Value = 99.9516
$ws.Range("F2").Value = 99.5684
$ws.Range("G2").Value = 99.7238
$ws.Range("H2").Value = 101.1178
$ws.Range("I2").Value = 102.3506
$ws.Range("A3").Value = "2014-11"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 99.5082
$ws.Range("D3").Value = 101.73
$ws.Range("E3").Value = 100.0501
$ws.Range("F3").Value = 99.6045
$ws.Range("G3").Value = 99.9002
$ws.Range("H3").Value = 101.1115
$ws.Range("I3").Value = 100.8773
$ws.Range("A4").Value = "2014-12"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 99.4953
$ws.Range("D4").Value = 101.3982
$ws.Range("E4").Value = 100.1108
$ws.Range("F4").Value = 99.44840000000001
$ws.Range("G4").Value = 99.8856
$ws.Range("H4").Value = 100.9919
$ws.Range("I4").Value = 100.6477
$ws.Range("A5").Value = "2014-01"
$ws.Range("B5").Value = 100.037
$ws.Range("C5").Value = 99.05370000000001
$ws.Range("D5").Value = 100.5211
$ws.Range("E5").Value = 100.5538
$ws.Range("F5").Value = 99.62260000000001
$ws.Range("G5").Value = 99.2272
$ws.Range("H5").Value = 99.261
$ws.Range("I5").Value = 99.2728
$ws.Range("A6").Value = "2014-02"
$ws.Range("B6").Value = 100.037
$ws.Range("C6").Value = 99.0081
$ws.Range("D6").Value = 100.2119
$ws.Range("E6").Value = 100.1572
$ws.Range("F6").Value = 99.8164
$ws.Range("G6").Value = 98.7713
$ws.Range("H6").Value = 99.23390000000001
$ws.Range("I6").Value = 99.2222
$ws.Range("A7").Value = "2014-03"
$ws.Range("B7").Value = 100.037
$ws.Range("C7").Value = 98.80329999999999
$ws.Range("D7").Value = 100.3398
$ws.Range("E7").Value = 99.90430000000001
$ws.Range("F7").Value = 100.4273
$ws.Range("G7").Value = 98.8805
$ws.Range("H7").Value = 99.3094
$ws.Range("I7").Value = 98.9126
$ws.Range("A8").Value = "2014-04"
$ws.Range("B8").Value = 100.037
$ws.Range("C8").Value = 98.8527
$ws.Range("D8").Value = 100.4112
$ws.Range("E8").Value = 99.9786
$ws.Range("F8").Value = 100.3923
$ws.Range("G8").Value = 99.5055
$ws.Range("H8").Value = 99.2055
$ws.Range("I8").Value = 99.027
$ws.Range("A9").Value = "2014-05"
$ws.Range("B9").Value = 100.037
$ws.Range("C9").Value = 98.78489999999999
$ws.Range("D9").Value = 100.4761
$ws.Range("E9").Value = 99.8557
$ws.Range("F9").Value = 100.6011
$ws.Range("G9").Value = 99.6082
$ws.Range("H9").Value = 99.86199999999999
$ws.Range("I9").Value = 99.4718
$ws.Range("A10").Value = "2014-06"
$ws.Range("B10").Value = 100.037
$ws.Range("C10").Value = 99.01300000000001
$ws.Range("D10").Value = 100.6774
$ws.Range("E10").Value = 99.73009999999999
$ws.Range("F10").Value = 100.6428
$ws.Range("G10").Value = 99.59229999999999
$ws.Range("H10").Value = 100.8062
$ws.Range("I10").Value = 99.7908
$ws.Range("A11").Value = "2014-07"
$ws.Range("B11").Value = 100.037
$ws.Range("C11").Value = 99.1572
$ws.Range("D11").Value = 101.6733
$ws.Range("E11").Value = 99.5065
$ws.Range("F11").Value = 100.5039
$ws.Range("G11").Value = 99.4867
$ws.Range("H11").Value = 100.2819
$ws.Range("I11").Value = 100.553
$ws.Range("A12").Value = "2014-08"
$ws.Range("B12").Value = 100.037
$ws.Range("C12").Value = 99.38679999999999
$ws.Range("D12").Value = 101.8727
$ws.Range("E12").Value = 99.7385
$ws.Range("F12").Value = 100.1936
$ws.Range("G12").Value = 99.81
$ws.Range("H12").Value = 100.3068
$ws.Range("I12").Value = 100.9798
$ws.Range("A13").Value = "2014-09"
$ws.Range("B13").Value = 100.037
$ws.Range("C13").Value = 99.3908
$ws.Range("D13").Value = 101.8108
$ws.Range("E13").Value = 99.2684
$ws.Range("F13").Value = 99.8105
$ws.Range("G13").Value = 99.73090000000001
$ws.Range("H13").Value = 100.4347
$ws.Range("I13").Value = 102.9681
$ws.Range("A14").Value = "2015-10"
$ws.Range("B14").Value = 100
$ws.Range("C14").Value = 99.2
$ws.Range("D14").Value = 96
$ws.Range("E14").Value = 99.5
$ws.Range("F14").Value = 99.90000000000001
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = 99.40000000000001
$ws.Range("I14").Value = 98.7
$ws.Range("A15").Value = "2015-11"
$ws.Range("B15").Value = 100
$ws.Range("C15").Value = 99.003
$ws.Range("D15").Value = 96.6386
$ws.Range("E15").Value = 99.4759
$ws.Range("F15").Value = 99.7824
$ws.Range("G15").Value = 99.8822
$ws.Range("H15").Value = 99.10939999999999
$ws.Range("I15").Value = 98.9177
$ws.Range("A16").Value = "2015-12"
$ws.Range("B16").Value = 100
$ws.Range("C16").Value = 98.91249999999999
$ws.Range("D16").Value = 96.90049999999999
$ws.Range("E16").Value = 99.3708
$ws.Range("F16").Value = 100.2661
$ws.Range("G16").Value = 99.7921
$ws.Range("H16").Value = 99.054
$ws.Range("I16").Value = 99.0728
$ws.Range("A17").Value = "2015-01"
$ws.Range("B17").Value = 100
$ws.Range("C17").Value = 99.4794
$ws.Range("D17").Value = 101.006
$ws.Range("E17").Value = 99.6387
$ws.Range("F17").Value = 99.00960000000001
$ws.Range("G17").Value = 100.2705
$ws.Range("H17").Value = 101.2659
$ws.Range("I17").Value = 100.7619
$ws.Range("A18").Value = "2015-02"
$ws.Range("B18").Value = 100
$ws.Range("C18").Value = 99.4432
$ws.Range("D18").Value = 97.1581
$ws.Range("E18").Value = 99.92910000000001
$ws.Range("F18").Value = 99.0094
$ws.Range("G18").Value = 100.5625
$ws.Range("H18").Value = 101.079
$ws.Range("I18").Value = 100.8341
$ws.Range("A19").Value = "2015-03"
$ws.Range("B19").Value = 100
$ws.Range("C19").Value = 99.5314
$ws.Range("D19").Value = 96.73260000000001
$ws.Range("E19").Value = 100.2477
$ws.Range("F19").Value = 98.3408
$ws.Range("G19").Value = 100.6334
$ws.Range("H19").Value = 100.5299
$ws.Range("I19").Value = 101.2043
$ws.Range("A20").Value = "2015-04"
$ws.Range("B20").Value = 100
$ws.Range("C20").Value = 99.3681
$ws.Range("D20").Value = 96.1495
$ws.Range("E20").Value = 100.2379
$ws.Range("F20").Value = 98.2757
$ws.Range("G20").Value = 100.4958
$ws.Range("H20").Value = 100.6857
$ws.Range("I20").Value = 100.6542
$ws.Range("A21").Value = "2015-05"
$ws.Range("B21").Value = 100
$ws.Range("C21").Value = 99.4357
$ws.Range("D21").Value = 96.19240000000001
$ws.Range("E21").Value = 100.2425
$ws.Range("F21").Value = 98.379
$ws.Range("G21").Value = 100.4675
$ws.Range("H21").Value = 100.4005
$ws.Range("I21").Value = 100.3292
$ws.Range("A22").Value = "2015-06"
$ws.Range("B22").Value = 100
$ws.Range("C22").Value = 99.6293
$ws.Range("D22").Value = 95.9738
$ws.Range("E22").Value = 100.2242
$ws.Range("F22").Value = 98.6105
$ws.Range("G22").Value = 100.1957
$ws.Range("H22").Value = 99.3927
$ws.Range("I22").Value = 100.0104
$ws.Range("A23").Value = "2015-07"
$ws.Range("B23").Value = 100
$ws.Range("C23").Value = 99.6323
$ws.Range("D23").Value = 95.1507
$ws.Range("E23").Value = 99.86239999999999
$ws.Range("F23").Value = 98.4776
$ws.Range("G23").Value = 99.99039999999999
$ws.Range("H23").Value = 99.8109
$ws.Range("I23").Value = 99.47709999999999
$ws.Range("A24").Value = "2015-08"
$ws.Range("B24").Value = 100
$ws.Range("C24").Value = 99.4447
$ws.Range("D24").Value = 95.6939
$ws.Range("E24").Value = 99.6909
$ws.Range("F24").Value = 99.2252
$ws.Range("G24").Value = 99.7139
$ws.Range("H24").Value = 99.7542
$ws.Range("I24").Value = 99.2824
$ws.Range("A25").Value = "2015-09"
$ws.Range("B25").Value = 100
$ws.Range("C25").Value = 99.3122
$ws.Range("D25").Value = 95.92619999999999
$ws.Range("E25").Value = 99.85899999999999
$ws.Range("F25").Value = 99.62820000000001
$ws.Range("G25").Value = 99.7841
$ws.Range("H25").Value = 99.6468
$ws.Range("I25").Value = 98.74120000000001
$ws.Range("A26").Value = "2016-10"
$ws.Range("B26").Value = 99.90000000000001
$ws.Range("C26").Value = 99.09999999999999
$ws.Range("D26").Value = 96.8
$ws.Range("E26").Value = 100.2
$ws.Range("F26").Value = 101.3
$ws.Range("G26").Value = 100.9
$ws.Range("H26").Value = 99.59999999999999
$ws.Range("I26").Value = 96.7
$ws.Range("A27").Value = "2016-11"
$ws.Range("B27").Value = 100
$ws.Range("C27").Value = 99.3
$ws.Range("D27").Value = 97.3
$ws.Range("E27").Value = 100.2
$ws.Range("F27").Value = 103.1
$ws.Range("G27").Value = 101.1
$ws.Range("H27").Value = 99.5
$ws.Range("I27").Value = 96.7
$ws.Range("A28").Value = "2016-12"
$ws.Range("B28").Value = 100.1
$ws.Range("C28").Value = 99.7
$ws.Range("D28").Value = 98.5
$ws.Range("E28").Value = 100.5
$ws.Range("F28").Value = 101.4
$ws.Range("G28").Value = 101.5
$ws.Range("H28").Value = 99.7
$ws.Range("I28").Value = 96.3
$ws.Range("A29").Value = "2016-01"
$ws.Range("B29").Value = 99.78019999999999
$ws.Range("C29").Value = 98.4186
$ws.Range("D29").Value = 93.2985
$ws.Range("E29").Value = 99.2244
$ws.Range("F29").Value = 100.3402
$ws.Range("G29").Value = 99.28830000000001
$ws.Range("H29").Value = 98.67610000000001
$ws.Range("I29").Value = 101.709
$ws.Range("A30").Value = "2016-02"
$ws.Range("B30").Value = 99.7881
$ws.Range("C30").Value = 98.4294
$ws.Range("D30").Value = 94.13760000000001
$ws.Range("E30").Value = 99.22499999999999
$ws.Range("F30").Value = 100.4768
$ws.Range("G30").Value = 99.1083
$ws.Range("H30").Value = 98.6584
$ws.Range("I30").Value = 102.289
$ws.Range("A31").Value = "2016-03"
$ws.Range("B31").Value = 99.9327
$ws.Range("C31").Value = 98.4764
$ws.Range("D31").Value = 93.10550000000001
$ws.Range("E31").Value = 99.0702
$ws.Range("F31").Value = 102.411
$ws.Range("G31").Value = 99.4195
$ws.Range("H31").Value = 99.1288
$ws.Range("I31").Value = 102.7353
$ws.Range("A32").Value = "2016-04"
$ws.Range("B32").Value = 99.95059999999999
$ws.Range("C32").Value = 98.4987
$ws.Range("D32").Value = 93.782
$ws.Range("E32").Value = 99.3749
$ws.Range("F32").Value = 102.7886
$ws.Range("G32").Value = 99.6947
$ws.Range("H32").Value = 98.9376
$ws.Range("I32").Value = 99.0698
$ws.Range("A33").Value = "2016-05"
$ws.Range("B33").Value = 100
$ws.Range("C33").Value = 98.7
$ws.Range("D33").Value = 94.40000000000001
$ws.Range("E33").Value = 99.5
$ws.Range("F33").Value = 102.7
$ws.Range("G33").Value = 100.1
$ws.Range("H33").Value = 99.2
$ws.Range("I33").Value = 98.8
$ws.Range("A34").Value = "2016-06"
$ws.Range("B34").Value = 100
$ws.Range("C34").Value = 98.90000000000001
$ws.Range("D34").Value = 95.40000000000001
$ws.Range("E34").Value = 99.59999999999999
$ws.Range("F34").Value = 102.6
$ws.Range("G34").Value = 100.6
$ws.Range("H34").Value = 99.2
$ws.Range("I34").Value = 97.3
$ws.Range("A35").Value = "2016-07"
$ws.Range("B35").Value = 100
$ws.Range("C35").Value = 98.90000000000001
$ws.Range("D35").Value = 95.5
$ws.Range("E35").Value = 99.90000000000001
$ws.Range("F35").Value = 103.2
$ws.Range("G35").Value = 100.7
$ws.Range("H35").Value = 99.3
$ws.Range("I35").Value = 96.90000000000001
$ws.Range("A36").Value = "2016-08"
$ws.Range("B36").Value = 100
$ws.Range("C36").Value = 98.7
$ws.Range("D36").Value = 96.40000000000001
$ws.Range("E36").Value = 100
$ws.Range("F36").Value = 101.4
$ws.Range("G36").Value = 100.6
$ws.Range("H36").Value = 99.40000000000001
$ws.Range("I36").Value = 97
$ws.Range("A37").Value = "2016-09"
$ws.Range("B37").Value = 100
$ws.Range("C37").Value = 99
$ws.Range("D37").Value = 96.40000000000001
$ws.Range("E37").Value = 100.2
$ws.Range("F37").Value = 102.2
$ws.Range("G37").Value = 100.8
$ws.Range("H37").Value = 99.5
$ws.Range("I37").Value = 96.8
$ws.Range("A38").Value = "2017-10"
$ws.Range("B38").Value = 94.8
$ws.Range("C38").Value = 101.3
$ws.Range("D38").Value = 100.7
$ws.Range("E38").Value = 101.3
$ws.Range("F38").Value = 102.3
$ws.Range("G38").Value = 102
$ws.Range("H38").Value = 100.4
$ws.Range("I38").Value = 99.3
$ws.Range("A39").Value = "2017-11"
$ws.Range("B39").Value = 94.8
$ws.Range("C39").Value = 101.5
$ws.Range("D39").Value = 100.8
$ws.Range("E39").Value = 101.5
$ws.Range("F39").Value = 102.6
$ws.Range("G39").Value = 101.8
$ws.Range("H39").Value = 100.3
$ws.Range("I39").Value = 99.2
$ws.Range("A40").Value = "2017-12"
$ws.Range("B40").Value = 94.90000000000001
$ws.Range("C40").Value = 101.4
$ws.Range("D40").Value = 100
$ws.Range("E40").Value = 101.6
$ws.Range("F40").Value = 102.5
$ws.Range("G40").Value = 102.1
$ws.Range("H40").Value = 100.4
$ws.Range("I40").Value = 99.59999999999999
$ws.Range("A41").Value = "2017-01"
$ws.Range("B41").Value = 96.90000000000001
$ws.Range("C41").Value = 100.4
$ws.Range("D41").Value = 98.8
$ws.Range("E41").Value = 101.1
$ws.Range("F41").Value = 102.3
$ws.Range("G41").Value = 101.9
$ws.Range("H41").Value = 100.2
$ws.Range("I41").Value = 96.2
$ws.Range("A42").Value = "2017-02"
$ws.Range("B42").Value = 97.2
$ws.Range("C42").Value = 100.5
$ws.Range("D42").Value = 98.90000000000001
$ws.Range("E42").Value = 100.8
$ws.Range("F42").Value = 101.1
$ws.Range("G42").Value = 102.3
$ws.Range("H42").Value = 100.1
$ws.Range("I42").Value = 96.2
$ws.Range("A43").Value = "2017-03"
$ws.Range("B43").Value = 97.2
$ws.Range("C43").Value = 100.7
$ws.Range("D43").Value = 101.1
$ws.Range("E43").Value = 101.2
$ws.Range("F43").Value = 102.2
$ws.Range("G43").Value = 102.5
$ws.Range("H43").Value = 100.3
$ws.Range("I43").Value = 97.09999999999999
$ws.Range("A44").Value = "2017-04"
$ws.Range("B44").Value = 97.2
$ws.Range("C44").Value = 100.8
$ws.Range("D44").Value = 100.5
$ws.Range("E44").Value = 101.1
$ws.Range("F44").Value = 102
$ws.Range("G44").Value = 102.3
$ws.Range("H44").Value = 100.7
$ws.Range("I44").Value = 99.8
$ws.Range("A45").Value = "2017-05"
$ws.Range("B45").Value = 97.2
$ws.Range("C45").Value = 100.8
$ws.Range("D45").Value = 99.59999999999999
$ws.Range("E45").Value = 101.1
$ws.Range("F45").Value = 103.6
$ws.Range("G45").Value = 101.9
$ws.Range("H45").Value = 100.7
$ws.Range("I45").Value = 99.7
$ws.Range("A46").Value = "2017-06"
$ws.Range("B46").Value = 97.2
$ws.Range("C46").Value = 100.7
$ws.Range("D46").Value = 98.8
$ws.Range("E46").Value = 101.2
$ws.Range("F46").Value = 102.4
$ws.Range("G46").Value = 101.9
$ws.Range("H46").Value = 100.7
$ws.Range("I46").Value = 100.4
$ws.Range("A47").Value = "2017-07"
$ws.Range("B47").Value = 97.2
$ws.Range("C47").Value = 100.8
$ws.Range("D47").Value = 99.7
$ws.Range("E47").Value = 100.9
$ws.Range("F47").Value = 104
$ws.Range("G47").Value = 101.9
$ws.Range("H47").Value = 100.5
$ws.Range("I47").Value = 100.1
$ws.Range("A48").Value = "2017-08"
$ws.Range("B48").Value = 97.2
$ws.Range("C48").Value = 101.1
$ws.Range("D48").Value = 100
$ws.Range("E48").Value = 100.8
$ws.Range("F48").Value = 103.1
$ws.Range("G48").Value = 101.9
$ws.Range("H48").Value = 100.7
$ws.Range("I48").Value = 99.40000000000001
$ws.Range("A49").Value = "2017-09"
$ws.Range("B49").Value = 97.09999999999999
$ws.Range("C49").Value = 101
$ws.Range("D49").Value = 100.5
$ws.Range("E49").Value = 101
$ws.Range("F49").Value = 102.1
$ws.Range("G49").Value = 101.5
$ws.Range("H49").Value = 100.6
$ws.Range("I49").Value = 98.90000000000001
